# NIT-9016844149.xlsx - "Estado de Cuenta" update
# - Remove the old/obsolete "Periodo Mora" (payment period) records (2402-2412),
#   keeping only periods through 2401.
# - Recalculate/refresh the totals shown for "Valor Mora" and "Cant. Periodos".
# - Re-apply the closing border/format to the new last data row of the table and
#   center the "Periodo Mora" column values, matching the refreshed table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table lives in rows 15 (header) through 38 (last detail row).
# Row 38 carries the special "closing" bottom border used only for the final
# row of the table. Before we drop the now-obsolete rows (27-38, periods
# 2402 through 2412, plus the duplicate 2401 row), copy that closing format
# onto row 26 (period 2401 / CC 71252818), which will become the new last row.
$ws.Range("B38:J38").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the obsolete trailing detail rows (old periods 2401-dup through 2412).
$ws.Rows("27:38").Delete()

# Center the "Periodo Mora" column for all remaining detail rows.
$ws.Range("E16:E26").HorizontalAlignment = -4108  # xlCenter

# Refresh the summary figures for the trimmed data set.
$ws.Range("E11").Value = 1116266
$ws.Range("F13").Value = 7
